$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) contain text that can look numeric
# (e.g. "0.670", "34.790.52"). Force text format first so Excel COM does
# not silently coerce these into Number/Date values, then restore the
# original "Normal" style so no stray number-format style is left on the cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "34.790.52"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "1.862.12"
$ws.Range("E3").Value = "  -2.66%  "
$ws.Range("E4").Value = "  -0.92%  "
$ws.Range("D5").Value = "244.03"
$ws.Range("E5").Value = "  -4.20%  "
$ws.Range("D6").Value = "0.670"
$ws.Range("E6").Value = "  -7.23%  "
$ws.Range("E7").Value = "  -1.00%  "
$ws.Range("D8").Value = "41.98"
$ws.Range("E8").Value = "  +2.83%  "
$ws.Range("D9").Value = "0.338"
$ws.Range("E9").Value = "  -6.07%  "
$ws.Range("D10").Value = "0.0729"
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  -2.50%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").Value = "12.79"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "2.127.69"
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").Value = "0.705"
$ws.Range("E14").Value = "  -2.42%  "
$ws.Range("D15").Value = "1.862.50"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("D16").Value = "4.79"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("D17").Value = "34.747.18"
$ws.Range("E17").Value = "  -1.72%  "
$ws.Range("D18").Value = "71.75"
$ws.Range("E18").Value = "  -3.91%  "
$ws.Range("D19").Value = "0.0₃0805"
$ws.Range("E19").Value = "  -6.01%  "
$ws.Range("D20").Value = "241.77"
$ws.Range("E20").Value = "  -1.18%  "
$ws.Range("D21").Value = "12.48"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("D22").Value = "4.84"
$ws.Range("E22").Value = "  -4.79%  "
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "2.48"
$ws.Range("E24").Value = "  +5.34%  "
$ws.Range("D25").Value = "2.15"
$ws.Range("E25").Value = "  -13.53%  "
$ws.Range("D26").Value = "162.73"
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("D27").Value = "8.28"
$ws.Range("E27").Value = "  -4.52%  "
$ws.Range("D28").Value = "17.96"
$ws.Range("E28").Value = "  -4.44%  "
$ws.Range("D29").Value = "0.125"
$ws.Range("E29").Value = "  -6.37%  "
$ws.Range("D31").Value = "1.71"
$ws.Range("E31").Value = "  +4.34%  "
$ws.Range("D32").Value = "4.13"
$ws.Range("E32").Value = "  -5.50%  "
$ws.Range("D33").Value = "0.0567"
$ws.Range("E33").Value = "  -2.89%  "
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").Value = "4.08"
$ws.Range("E35").Value = "  -3.55%  "
$ws.Range("D36").Value = "0.820"
$ws.Range("E36").Value = "  -10.99%  "
$ws.Range("D37").Value = "1.92"
$ws.Range("E37").Value = "  -5.77%  "
$ws.Range("D38").Value = "1.50"
$ws.Range("E38").Value = "  -24.71%  "
$ws.Range("D39").Value = "97.12"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "16.88"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "0.0660"
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("D42").Value = "0.0209"
$ws.Range("E42").Value = "  -5.05%  "
$ws.Range("D43").Value = "1.06"
$ws.Range("E43").Value = "  -5.41%  "
$ws.Range("D44").Value = "0.0841"
$ws.Range("E44").Value = "  +13.62%  "
$ws.Range("D45").Value = "1.274.62"
$ws.Range("E45").Value = "  -5.05%  "
$ws.Range("D46").Value = "2.28"
$ws.Range("E46").Value = "  -6.54%  "
$ws.Range("E48").Value = "  -1.71%  "
$ws.Range("D49").Value = "11.85"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").Value = "6.22"
$ws.Range("E50").Value = "  -8.23%  "
$ws.Range("D51").Value = "42.09"
$ws.Range("E51").Value = "  -7.15%  "

$ws.Range("D2:E51").Style = "Normal"
